$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11529.429
$ws.Range("I70").Value = 12925.375
$ws.Range("J70").Value = 9668.166999999999
$ws.Range("K70").Value = 38776.125
$ws.Range("L70").Value = 29004.501
$ws.Range("M70").Value = -38506.125
$ws.Range("N70").Value = -29544.501
$ws.Range("H73").Value = 11529.429
$ws.Range("I73").Value = 12925.375
$ws.Range("J73").Value = 9668.166999999999
$ws.Range("K73").Value = 38776.125
$ws.Range("L73").Value = 29004.501
$ws.Range("M73").Value = -37840.125
$ws.Range("N73").Value = -30876.501
$ws.Range("H125").Value = 142858740
$ws.Range("I125").Value = 250001000
$ws.Range("J125").Value = 2366.3333
$ws.Range("K125").Value = 2250009000
$ws.Range("L125").Value = 21296.9997
$ws.Range("M125").Value = -2250006540
$ws.Range("H135").Value = 1818944.9
$ws.Range("I135").Value = 2222706.5
$ws.Range("J135").Value = 2018
$ws.Range("K135").Value = 20004358.5
$ws.Range("L135").Value = 18162
$ws.Range("M135").Value = -20001823.5
$ws.Range("N135").Value = -23232
$ws.Range("H137").Value = 4274.7617
$ws.Range("I137").Value = 2263.6924
$ws.Range("J137").Value = 7542.75
$ws.Range("K137").Value = 6791.0772
$ws.Range("L137").Value = 22628.25
$ws.Range("M137").Value = -4241.0772
$ws.Range("N137").Value = -27728.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3642049.8
$ws.Range("I32").Value = 3777787.5
$ws.Range("J32").Value = 44999.5
$ws.Range("K32").Value = 3777787.5
$ws.Range("L32").Value = 44999.5
$ws.Range("M32").Value = -3777500.5
$ws.Range("H61").Value = 6938.391
$ws.Range("I61").Value = 2158.2144
$ws.Range("J61").Value = 14374.223
$ws.Range("K61").Value = 2158.2144
$ws.Range("L61").Value = 14374.223
$ws.Range("M61").Value = -1946.2144
$ws.Range("N61").Value = -14798.223
$ws.Range("H74").Value = 41356.96
$ws.Range("I74").Value = 57515.61
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 57515.61
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -56641.61
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 41356.96
$ws.Range("I77").Value = 57515.61
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 287578.05
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -283210.05
$ws.Range("N77").Value = -33736
$ws.Range("H110").Value = 13338453
$ws.Range("I110").Value = 6525.0527
$ws.Range("J110").Value = 55556224
$ws.Range("K110").Value = 6525.0527
$ws.Range("L110").Value = 55556224
$ws.Range("M110").Value = -4480.0527
$ws.Range("H136").Value = 6938.391
$ws.Range("I136").Value = 2158.2144
$ws.Range("J136").Value = 14374.223
$ws.Range("K136").Value = 6474.6432
$ws.Range("L136").Value = 43122.669
$ws.Range("M136").Value = -3924.6432
$ws.Range("N136").Value = -48222.669

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6471.5586
$ws.Range("I134").Value = 2866.65
$ws.Range("J134").Value = 11621.429
$ws.Range("K134").Value = 8599.950000000001
$ws.Range("L134").Value = 34864.287
$ws.Range("M134").Value = -6064.950000000001
$ws.Range("N134").Value = -39934.287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6355
$ws.Range("I31").Value = 1766.625
$ws.Range("J31").Value = 10759.84
$ws.Range("K31").Value = 1766.625
$ws.Range("L31").Value = 10759.84
$ws.Range("M31").Value = -1471.625
$ws.Range("N31").Value = -11349.84
$ws.Range("H34").Value = 6355
$ws.Range("I34").Value = 1766.625
$ws.Range("J34").Value = 10759.84
$ws.Range("K34").Value = 1766.625
$ws.Range("L34").Value = 10759.84
$ws.Range("M34").Value = -1564.625
$ws.Range("N34").Value = -11163.84
$ws.Range("H62").Value = 32412406
$ws.Range("I62").Value = 48614892
$ws.Range("J62").Value = 7431.3335
$ws.Range("K62").Value = 48614892
$ws.Range("L62").Value = 7431.3335
$ws.Range("M62").Value = -48614268
$ws.Range("N62").Value = -8679.333500000001
$ws.Range("H65").Value = 32412406
$ws.Range("I65").Value = 48614892
$ws.Range("J65").Value = 7431.3335
$ws.Range("K65").Value = 243074460
$ws.Range("L65").Value = 37156.6675
$ws.Range("M65").Value = -243071340
$ws.Range("N65").Value = -43396.6675
$ws.Range("H99").Value = 4962.846
$ws.Range("I99").Value = 3500.8
$ws.Range("J99").Value = 5876.625
$ws.Range("K99").Value = 3500.8
$ws.Range("L99").Value = 5876.625
$ws.Range("M99").Value = -2002.8
$ws.Range("H126").Value = 4962.846
$ws.Range("I126").Value = 3500.8
$ws.Range("J126").Value = 5876.625
$ws.Range("K126").Value = 10502.4
$ws.Range("L126").Value = 17629.875
$ws.Range("M126").Value = -8032.400000000001
$ws.Range("H134").Value = 5788.162
$ws.Range("I134").Value = 2277.875
$ws.Range("J134").Value = 8462.666999999999
$ws.Range("K134").Value = 6833.625
$ws.Range("L134").Value = 25388.001
$ws.Range("M134").Value = -4298.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1869.6316
$ws.Range("I5").Value = 884.9167
$ws.Range("J5").Value = 3557.7144
$ws.Range("K5").Value = 2654.7501
$ws.Range("L5").Value = 10673.1432
$ws.Range("M5").Value = -2542.7501
$ws.Range("N5").Value = -10897.1432
$ws.Range("H22").Value = 704
$ws.Range("I22").Value = 140
$ws.Range("J22").Value = 986
$ws.Range("K22").Value = 420
$ws.Range("L22").Value = 2958
$ws.Range("M22").Value = -251
$ws.Range("N22").Value = -3296
$ws.Range("H27").Value = 704
$ws.Range("I27").Value = 140
$ws.Range("J27").Value = 986
$ws.Range("K27").Value = 420
$ws.Range("L27").Value = 2958
$ws.Range("M27").Value = -318
$ws.Range("N27").Value = -3162
$ws.Range("H68").Value = 3142.818
$ws.Range("I68").Value = 659.8
$ws.Range("J68").Value = 5212
$ws.Range("K68").Value = 1979.4
$ws.Range("L68").Value = 15636
$ws.Range("M68").Value = -1168.4
$ws.Range("N68").Value = -17258
$ws.Range("H71").Value = 3142.818
$ws.Range("I71").Value = 659.8
$ws.Range("J71").Value = 5212
$ws.Range("K71").Value = 5938.2
$ws.Range("L71").Value = 46908
$ws.Range("M71").Value = -1882.2
$ws.Range("N71").Value = -55020
$ws.Range("H87").Value = 790.75
$ws.Range("I87").Value = 790.75
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 2372.25
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -1124.25
$ws.Range("H90").Value = 790.75
$ws.Range("I90").Value = 790.75
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 7116.75
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -876.75
$ws.Range("H135").Value = 1869.6316
$ws.Range("I135").Value = 884.9167
$ws.Range("J135").Value = 3557.7144
$ws.Range("K135").Value = 7964.2503
$ws.Range("L135").Value = 32019.4296
$ws.Range("M135").Value = -5429.2503
$ws.Range("N135").Value = -37089.4296
$ws.Range("H137").Value = 156150.61
$ws.Range("I137").Value = 112967.78
$ws.Range("J137").Value = 253312
$ws.Range("K137").Value = 338903.34
$ws.Range("L137").Value = 759936
$ws.Range("M137").Value = -333803.34

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2446.6667
$ws.Range("I80").Value = 2465.375
$ws.Range("J80").Value = 2297
$ws.Range("K80").Value = 2465.375
$ws.Range("L80").Value = 2297
$ws.Range("M80").Value = -1467.375
$ws.Range("H83").Value = 2446.6667
$ws.Range("I83").Value = 2465.375
$ws.Range("J83").Value = 2297
$ws.Range("K83").Value = 12326.875
$ws.Range("L83").Value = 11485
$ws.Range("M83").Value = -7334.875
$ws.Range("H113").Value = 8979.791999999999
$ws.Range("I113").Value = 6987.375
$ws.Range("J113").Value = 9976
$ws.Range("K113").Value = 6987.375
$ws.Range("L113").Value = 9976
$ws.Range("M113").Value = -4817.375
$ws.Range("H132").Value = 6207.5293
$ws.Range("I132").Value = 2155.1
$ws.Range("J132").Value = 11996.714
$ws.Range("K132").Value = 6465.299999999999
$ws.Range("L132").Value = 35990.142
$ws.Range("M132").Value = -3935.299999999999
$ws.Range("N132").Value = -41050.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5119.9287
$ws.Range("I7").Value = 4172.5
$ws.Range("J7").Value = 6383.1665
$ws.Range("K7").Value = 4172.5
$ws.Range("L7").Value = 6383.1665
$ws.Range("M7").Value = -4060.5
$ws.Range("N7").Value = -6607.1665
$ws.Range("H40").Value = 5538.737
$ws.Range("I40").Value = 2472.5557
$ws.Range("J40").Value = 8298.299999999999
$ws.Range("K40").Value = 2472.5557
$ws.Range("L40").Value = 8298.299999999999
$ws.Range("M40").Value = -2336.5557
$ws.Range("H43").Value = 9142.857
$ws.Range("I43").Value = 7000
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 7000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10386
$ws.Range("M43").Value = -6807
$ws.Range("H68").Value = 4807.9473
$ws.Range("I68").Value = 2805.7778
$ws.Range("J68").Value = 6609.9
$ws.Range("K68").Value = 2805.7778
$ws.Range("L68").Value = 6609.9
$ws.Range("M68").Value = -2056.7778
$ws.Range("N68").Value = -8107.9
$ws.Range("H71").Value = 4807.9473
$ws.Range("I71").Value = 2805.7778
$ws.Range("J71").Value = 6609.9
$ws.Range("K71").Value = 14028.889
$ws.Range("L71").Value = 33049.5
$ws.Range("M71").Value = -10284.889
$ws.Range("N71").Value = -40537.5
$ws.Range("H93").Value = 812.875
$ws.Range("I93").Value = 583.8333
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 583.8333
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 664.1667
$ws.Range("H122").Value = 4005.2593
$ws.Range("I122").Value = 2667.111
$ws.Range("J122").Value = 6681.5557
$ws.Range("K122").Value = 8001.333
$ws.Range("L122").Value = 20044.6671
$ws.Range("M122").Value = -5551.333
$ws.Range("N122").Value = -24944.6671
$ws.Range("H126").Value = 5119.9287
$ws.Range("I126").Value = 4172.5
$ws.Range("J126").Value = 6383.1665
$ws.Range("K126").Value = 12517.5
$ws.Range("L126").Value = 19149.4995
$ws.Range("M126").Value = -10047.5
$ws.Range("N126").Value = -24089.4995
$ws.Range("H132").Value = 6615.6665
$ws.Range("I132").Value = 3226.6316
$ws.Range("J132").Value = 9092.27
$ws.Range("K132").Value = 9679.8948
$ws.Range("L132").Value = 27276.81
$ws.Range("M132").Value = -7149.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6249.143
$ws.Range("I62").Value = 6848.8
$ws.Range("J62").Value = 4750
$ws.Range("K62").Value = 6848.8
$ws.Range("L62").Value = 4750
$ws.Range("M62").Value = -6224.8
$ws.Range("N62").Value = -5998
$ws.Range("H65").Value = 6249.143
$ws.Range("I65").Value = 6848.8
$ws.Range("J65").Value = 4750
$ws.Range("K65").Value = 34244
$ws.Range("L65").Value = 23750
$ws.Range("M65").Value = -31124
$ws.Range("N65").Value = -29990
$ws.Range("H81").Value = 22320732
$ws.Range("I81").Value = 126649.25
$ws.Range("J81").Value = 40076000
$ws.Range("K81").Value = 253298.5
$ws.Range("L81").Value = 80152000
$ws.Range("M81").Value = -252237.5
$ws.Range("N81").Value = -80154122
$ws.Range("H84").Value = 22320732
$ws.Range("I84").Value = 126649.25
$ws.Range("J84").Value = 40076000
$ws.Range("K84").Value = 1266492.5
$ws.Range("L84").Value = 400760000
$ws.Range("M84").Value = -1261188.5
$ws.Range("N84").Value = -400770608
$ws.Range("H122").Value = 158218.92
$ws.Range("I122").Value = 213916.16
$ws.Range("J122").Value = 7040.7144
$ws.Range("K122").Value = 641748.48
$ws.Range("L122").Value = 21122.1432
$ws.Range("M122").Value = -639298.48
$ws.Range("H126").Value = 981.6667
$ws.Range("I126").Value = 981.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2945.0001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -475.0001000000002
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 9104.647000000001
$ws.Range("I132").Value = 14997.714
$ws.Range("J132").Value = 4979.5
$ws.Range("K132").Value = 44993.142
$ws.Range("L132").Value = 14938.5
$ws.Range("M132").Value = -42463.142

Write-Host "done"